# Apply the target edit to the presentation:
#   1. Re-colour the presentation's (slide master's) theme colour scheme
#      from the "Red Violet" swatch used by the "Integral" theme to the
#      standard "Office" swatch (ppt/theme/theme2.xml : a:clrScheme).
#   2. Re-point the three data tables (on the three consecutive summary
#      slides) at the new built-in table style
#      {7B572EDB-FFE4-46F9-89C4-7BF45F019F18} instead of the old
#      {9DCBE013-7E29-4F01-A294-C6E94D7113CD}.

$p = $ppt.ActivePresentation

# --- 1. Recolour the theme (12 theme colour slots, in clrScheme order:
#        dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) -----------------
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme
$tcs.Colors(1).RGB  = 0         # dk1      000000
$tcs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388   # dk2      44546A
$tcs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407     # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308  # accent5  4472C4
$tcs.Colors(10).RGB = 4697456   # accent6  70AD47
$tcs.Colors(11).RGB = 12673797  # hlink    0563C1
$tcs.Colors(12).RGB = 7491477   # folHlink 954F72

# --- 2. Swap the table style used by the three summary tables ----------
$oldStyleId = "{9DCBE013-7E29-4F01-A294-C6E94D7113CD}"
$newStyleId = "{7B572EDB-FFE4-46F9-89C4-7BF45F019F18}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
